$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("A1").Value = "Factura"

# Row 2
$ws.Range("A2").Value = "TEMP0000049"
$ws.Range("C2").Value = "27/11/2024 17:15:00"
$ws.Range("D2").Value = "27/11/2024 17:17:00"
$ws.Range("G2").Value = "00:02"

# Row 3
$ws.Range("A3").Value = "TEMP0000050"
$ws.Range("B3").Value = "BBB001"
$ws.Range("C3").Value = "27/11/2024 17:16:00"
$ws.Range("D3").Value = "27/11/2024 17:17:00"
$ws.Range("E3").Value = "Carro"
$ws.Range("F3").Value = 4000
$ws.Range("G3").Value = "00:01"
$ws.Range("H3").Value = 4000

# Row 4
$ws.Range("A4").Value = "TEMP0000051"
$ws.Range("B4").Value = "CCC01"
$ws.Range("C4").Value = "28/11/2024 09:45:00"
$ws.Range("D4").Value = "28/11/2024 09:46:00"
$ws.Range("G4").Value = "00:01"

# Row 5
$ws.Range("A5").Value = "TEMP0000052"
$ws.Range("B5").Value = "BBB002"
$ws.Range("C5").Value = "28/11/2024 09:45:00"
$ws.Range("D5").Value = "28/11/2024 09:46:00"
$ws.Range("E5").Value = "Carro"
$ws.Range("F5").Value = 4000
$ws.Range("G5").Value = "00:01"
$ws.Range("H5").Value = 4000

# Row 6
$ws.Range("A6").Value = "TEMP0000053"
$ws.Range("B6").Value = "DDD01"
$ws.Range("C6").Value = "28/11/2024 09:47:00"
$ws.Range("D6").Value = "28/11/2024 09:47:00"
$ws.Range("G6").Value = "00:00"

# Row 7
$ws.Range("A7").Value = "TEMP0000054"
$ws.Range("B7").Value = "AAA03"
$ws.Range("C7").Value = "28/11/2024 10:00:00"
$ws.Range("D7").Value = "28/11/2024 10:01:00"
$ws.Range("G7").Value = "00:01"

# Row 8
$ws.Range("A8").Value = "TEMP0000055"
$ws.Range("B8").Value = "EEE003"
$ws.Range("C8").Value = "28/11/2024 10:00:00"
$ws.Range("D8").Value = "28/11/2024 10:02:00"
$ws.Range("E8").Value = "Carro"
$ws.Range("F8").Value = 4000
$ws.Range("G8").Value = "00:02"
$ws.Range("H8").Value = 4000

# Row 9
$ws.Range("A9").Value = "TEMP0000056"
$ws.Range("B9").Value = "FFF01"
$ws.Range("C9").Value = "28/11/2024 11:19:00"
$ws.Range("D9").Value = "28/11/2024 11:20:00"
$ws.Range("F9").Value = 1500
$ws.Range("G9").Value = "00:01"
$ws.Range("H9").Value = 1500
